$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $newValue)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextValue "D2" "60.297.86"
Set-TextValue "E2" "  +3.99%  "
Set-TextValue "D3" "2.339.84"
Set-TextValue "E3" "  +2.41%  "
Set-TextValue "E4" "  +0.06%  "
Set-TextValue "D5" "545.37"
Set-TextValue "E5" "  +2.89%  "
Set-TextValue "D6" "131.82"
Set-TextValue "E6" "  +0.85%  "
Set-TextValue "E7" "  +0.02%  "
Set-TextValue "E8" "  +0.60%  "
Set-TextValue "D9" "2.337.07"
Set-TextValue "E9" "  +2.44%  "
Set-TextValue "E10" "  +1.88%  "
Set-TextValue "E11" "  +1.12%  "
Set-TextValue "E12" "  +0.88%  "
Set-TextValue "D13" "0.334"
Set-TextValue "E13" "  +1.79%  "
Set-TextValue "D14" "23.83"
Set-TextValue "E14" "  +2.07%  "
Set-TextValue "D15" "2.755.40"
Set-TextValue "E15" "  +2.43%  "
Set-TextValue "D16" "60.257.80"
Set-TextValue "E16" "  +4.14%  "
Set-TextValue "E17" "  +1.37%  "
Set-TextValue "D18" "2.341.95"
Set-TextValue "E18" "  +2.40%  "
Set-TextValue "E19" "  +1.07%  "
Set-TextValue "D20" "4.16"
Set-TextValue "E20" "  -0.01%  "
Set-TextValue "D21" "6.79"
Set-TextValue "E21" "  +6.39%  "
Set-TextValue "D22" "313.88"
Set-TextValue "E22" "  +0.96%  "
Set-TextValue "E23" "  -0.22%  "
Set-TextValue "D24" "63.51"
Set-TextValue "E24" "  +2.02%  "
Set-TextValue "D25" "0.172"
Set-TextValue "E25" "  +2.89%  "
Set-TextValue "E26" "  +0.02%  "
Set-TextValue "D27" "7.92"
Set-TextValue "E27" "  -0.63%  "
Set-TextValue "E28" "  +8.46%  "
Set-TextValue "E29" "  +2.86%  "
Set-TextValue "D30" "171.73"
Set-TextValue "E31" "  +13.79%  "
Set-TextValue "D32" "0.0₃0731"
Set-TextValue "E32" "  +2.42%  "
Set-TextValue "D33" "5.94"
Set-TextValue "E33" "  +4.12%  "
Set-TextValue "E34" "  +13.28%  "
Set-TextValue "D35" "0.382"
Set-TextValue "E35" "  +1.74%  "
Set-TextValue "D36" "18.02"
Set-TextValue "E36" "  +1.69%  "
Set-TextValue "E38" "  +0.01%  "
Set-TextValue "E39" "  +7.25%  "
Set-TextValue "D40" "321.38"
Set-TextValue "E40" "  +12.27%  "
Set-TextValue "D41" "38.15"
Set-TextValue "E41" "  -0.66%  "
Set-TextValue "E42" "  +3.06%  "
Set-TextValue "D43" "140.95"
Set-TextValue "E43" "  +0.36%  "
Set-TextValue "E44" "  +1.78%  "
Set-TextValue "D45" "0.0946"
Set-TextValue "E45" "  -0.04%  "
Set-TextValue "E46" "  +8.95%  "
Set-TextValue "D47" "0.0497"
Set-TextValue "E47" "  +0.69%  "
Set-TextValue "E48" "  +1.65%  "
Set-TextValue "E49" "  +2.18%  "
Set-TextValue "B50" "BabyDogeCoin"
Set-TextValue "C50" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D50" "0.0₆0212"
Set-TextValue "E50" "  +18.69%  "
Set-TextValue "B51" "WhiteBITCoin"
Set-TextValue "C51" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D51" "11.03"
Set-TextValue "E51" "  +0.84%  "

Write-Host "Applied" 80 "cell updates."
